# Weekly update: a new daily price record is inserted for
# "Vega Central Mapocho de Santiago - Bruselas (repollito)".
# The new record belongs chronologically before the existing row 54,
# so a new row is inserted at position 54 (pushing the former rows
# 54-88 down to 55-89) and populated with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 54; everything from the old row 54
# downward shifts down by one (old row 54 -> 55, ..., old row 88 -> 89).
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the new observation.
$ws.Cells.Item(54, 1).Value = 9
$ws.Cells.Item(54, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(54, 3).Value = "Metropolitana"
$ws.Cells.Item(54, 4).Value = 45086
$ws.Cells.Item(54, 5).Value = 13
$ws.Cells.Item(54, 6).Value = 100112035
$ws.Cells.Item(54, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(54, 8).Value = "Sin especificar"
$ws.Cells.Item(54, 9).Value = "Primera"
$ws.Cells.Item(54, 10).Value = 70
$ws.Cells.Item(54, 11).Value = 16000
$ws.Cells.Item(54, 12).Value = 18000
$ws.Cells.Item(54, 13).Value = 17000
$ws.Cells.Item(54, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(54, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(54, 16).Value = 1133
$ws.Cells.Item(54, 17).Value = 15
$ws.Cells.Item(54, 18).Value = "Hortaliza"
